$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 (old observation, 2022-03-08 / volumen 60) is duplicated down into
# new row 7 before row 6 itself is updated to the later weekly observation
# (2022-03-17 / volumen 100).
$rowSrc = 6
$rowNew = 7
$lastCol = 18

for ($col = 1; $col -le $lastCol; $col++) {
    $srcCell = $ws.Cells.Item($rowSrc, $col)
    $dstCell = $ws.Cells.Item($rowNew, $col)
    $dstCell.Value = $srcCell.Value2
}

# Column D (Fecha) carries the date number format style; copy it across too.
$ws.Cells.Item($rowNew, 4).NumberFormat = $ws.Cells.Item($rowSrc, 4).NumberFormat

# Update row 6 with the new weekly observation: later date, higher volume.
$ws.Cells.Item($rowSrc, 4).Value = 44637
$ws.Cells.Item($rowSrc, 10).Value = 100
